$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row permutation mapping: new row N gets the old content of row mapping[N]
# (columns D, I, J, K, L, M, N, P, Q; other columns are identical across rows)
$mapping = @{2 = 15; 3 = 8; 4 = 9; 5 = 14; 6 = 10; 7 = 3; 8 = 4; 9 = 2; 10 = 11; 11 = 18; 12 = 6; 13 = 7; 14 = 13; 15 = 12; 16 = 17; 17 = 16; 18 = 5}

$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")

# Snapshot the "before" values for each tracked column/row since we will
# be overwriting rows in place and some rows feed others as sources.
$snapshot = @{}
for ($r = 2; $r -le 18; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}

Write-Output "done"
